# Since we send two images, we need twice the number of tokens.
# The "Token Estimation / request" formulas in E35, E37 and E39 are
# doubled to account for sending two images per request.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E35").Formula = "=((D23/D26*E23/E26) +2 +D24)*2"
$ws.Range("E37").Formula = "=2*((D23/D27*E23/E27)+2+D25)"
$ws.Range("E39").Formula = "=2*((D23/D27*E23/E27)+2+D25)"

$ws.Range("E40").Select() | Out-Null
